$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4337.967
$ws.Range("I40").Value = 2171.2
$ws.Range("K40").Value = 2171.2
$ws.Range("M40").Value = -1996.2

$ws.Range("H64").Value = 8078.1665
$ws.Range("I64").Value = 5691
$ws.Range("J64").Value = 8248.679
$ws.Range("K64").Value = 5691
$ws.Range("L64").Value = 8248.679
$ws.Range("N64").Value = -8744.679
$ws.Range("M64").Value = -5443

$ws.Range("H67").Value = 8078.1665
$ws.Range("I67").Value = 5691
$ws.Range("J67").Value = 8248.679
$ws.Range("K67").Value = 5691
$ws.Range("L67").Value = 8248.679
$ws.Range("N67").Value = -9964.679
$ws.Range("M67").Value = -4833

$ws.Range("H92").Value = 2909
$ws.Range("I92").Value = 2171.6667
$ws.Range("J92").Value = 7333
$ws.Range("K92").Value = 2171.6667
$ws.Range("L92").Value = 7333
$ws.Range("M92").Value = -923.6667000000002
$ws.Range("N92").Value = -9829

$ws.Range("H112").Value = 7750.5625
$ws.Range("J112").Value = 8200.6
$ws.Range("L112").Value = 24601.8
$ws.Range("N112").Value = -26817.8

$ws.Range("H129").Value = 15153111
$ws.Range("I129").Value = 26316512
$ws.Range("J129").Value = 2780.0715
$ws.Range("K129").Value = 78949536
$ws.Range("L129").Value = 8340.2145
$ws.Range("M129").Value = -78944536
$ws.Range("N129").Value = -18340.2145

$ws.Range("H132").Value = 3028.196
$ws.Range("I132").Value = 3342.077
$ws.Range("J132").Value = 2008.0834
$ws.Range("K132").Value = 10026.231
$ws.Range("L132").Value = 6024.2502
$ws.Range("M132").Value = -7496.231
$ws.Range("N132").Value = -11084.2502

$ws.Range("H138").Value = 2650.6167
$ws.Range("I138").Value = 1328.2593
$ws.Range("J138").Value = 3732.5454
$ws.Range("K138").Value = 3984.7779
$ws.Range("L138").Value = 11197.6362
$ws.Range("M138").Value = 1155.2221
$ws.Range("N138").Value = -21477.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1178940
$ws.Range("I2").Value = 1379818
$ws.Range("K2").Value = 1379818
$ws.Range("M2").Value = -1379705

$ws.Range("H32").Value = 11326.057
$ws.Range("I32").Value = 6247.9185
$ws.Range("J32").Value = 22636.455
$ws.Range("K32").Value = 6247.9185
$ws.Range("L32").Value = 22636.455
$ws.Range("M32").Value = -5960.9185
$ws.Range("N32").Value = -23210.455

$ws.Range("H61").Value = 4958.1816
$ws.Range("I61").Value = 4949.778
$ws.Range("J61").Value = 4996
$ws.Range("K61").Value = 4949.778
$ws.Range("L61").Value = 4996
$ws.Range("M61").Value = -4737.778
$ws.Range("N61").Value = -5420

$ws.Range("H116").Value = 1178940
$ws.Range("I116").Value = 1379818
$ws.Range("K116").Value = 1379818
$ws.Range("M116").Value = -1377524

$ws.Range("H122").Value = 2627734.5
$ws.Range("I122").Value = 3292014.8
$ws.Range("J122").Value = 1742027.6
$ws.Range("K122").Value = 9876044.399999999
$ws.Range("L122").Value = 5226082.800000001
$ws.Range("M122").Value = -9873594.399999999
$ws.Range("N122").Value = -5230982.800000001

$ws.Range("H136").Value = 4958.1816
$ws.Range("I136").Value = 4949.778
$ws.Range("J136").Value = 4996
$ws.Range("K136").Value = 14849.334
$ws.Range("L136").Value = 14988
$ws.Range("M136").Value = -12299.334
$ws.Range("N136").Value = -20088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1178940
$ws.Range("I3").Value = 1379818
$ws.Range("K3").Value = 1379818
$ws.Range("M3").Value = -1379704

$ws.Range("H20").Value = 13890093
$ws.Range("I20").Value = 18519522
$ws.Range("J20").Value = 1804.3334
$ws.Range("K20").Value = 18519522
$ws.Range("L20").Value = 1804.3334
$ws.Range("M20").Value = -18519275
$ws.Range("N20").Value = -2298.3334

$ws.Range("H64").Value = 1191.4615
$ws.Range("I64").Value = 497.33334
$ws.Range("J64").Value = 1399.7
$ws.Range("K64").Value = 497.33334
$ws.Range("L64").Value = 1399.7
$ws.Range("M64").Value = -272.33334
$ws.Range("N64").Value = -1849.7

$ws.Range("H67").Value = 1191.4615
$ws.Range("I67").Value = 497.33334
$ws.Range("J67").Value = 1399.7
$ws.Range("K67").Value = 497.33334
$ws.Range("L67").Value = 1399.7
$ws.Range("M67").Value = 282.66666
$ws.Range("N67").Value = -2959.7

$ws.Range("H86").Value = 14446355
$ws.Range("I86").Value = 18842192
$ws.Range("J86").Value = 2892.7144
$ws.Range("K86").Value = 18842192
$ws.Range("L86").Value = 2892.7144
$ws.Range("M86").Value = -18841069
$ws.Range("N86").Value = -5138.7144

$ws.Range("H89").Value = 14446355
$ws.Range("I89").Value = 18842192
$ws.Range("J89").Value = 2892.7144
$ws.Range("K89").Value = 94210960
$ws.Range("L89").Value = 14463.572
$ws.Range("M89").Value = -94205344
$ws.Range("N89").Value = -25695.572

$ws.Range("H105").Value = 3907623
$ws.Range("I105").Value = 4167931.2
$ws.Range("K105").Value = 4167931.2
$ws.Range("M105").Value = -4166184.2

$ws.Range("H134").Value = 4470.6665
$ws.Range("I134").Value = 1122.138
$ws.Range("K134").Value = 3366.414
$ws.Range("M134").Value = -831.4139999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1505.2222
$ws.Range("I16").Value = 1382.5385
$ws.Range("J16").Value = 1824.2
$ws.Range("K16").Value = 1382.5385
$ws.Range("L16").Value = 1824.2
$ws.Range("M16").Value = -1095.5385
$ws.Range("N16").Value = -2398.2

$ws.Range("H31").Value = 19368.39
$ws.Range("I31").Value = 2811.3125
$ws.Range("J31").Value = 25529.162
$ws.Range("K31").Value = 2811.3125
$ws.Range("L31").Value = 25529.162
$ws.Range("M31").Value = -2516.3125
$ws.Range("N31").Value = -26119.162

$ws.Range("H34").Value = 19368.39
$ws.Range("I34").Value = 2811.3125
$ws.Range("J34").Value = 25529.162
$ws.Range("K34").Value = 2811.3125
$ws.Range("L34").Value = 25529.162
$ws.Range("M34").Value = -2609.3125
$ws.Range("N34").Value = -25933.162

$ws.Range("H58").Value = 6562.1304
$ws.Range("I58").Value = 7571.0625
$ws.Range("K58").Value = 7571.0625
$ws.Range("M58").Value = -7368.0625

$ws.Range("H86").Value = 9075.65
$ws.Range("I86").Value = 7512.231
$ws.Range("K86").Value = 7512.231
$ws.Range("M86").Value = -6389.231

$ws.Range("H89").Value = 9075.65
$ws.Range("I89").Value = 7512.231
$ws.Range("K89").Value = 37561.155
$ws.Range("M89").Value = -31945.155

$ws.Range("H113").Value = 1505.2222
$ws.Range("I113").Value = 1382.5385
$ws.Range("J113").Value = 1824.2
$ws.Range("K113").Value = 1382.5385
$ws.Range("L113").Value = 1824.2
$ws.Range("M113").Value = 787.4614999999999
$ws.Range("N113").Value = -6164.2

$ws.Range("H122").Value = 1986.55
$ws.Range("I122").Value = 1842.0769
$ws.Range("K122").Value = 5526.2307
$ws.Range("M122").Value = -3076.2307

$ws.Range("H132").Value = 120020.94
$ws.Range("I132").Value = 78573.234
$ws.Range("K132").Value = 235719.702
$ws.Range("M132").Value = -233189.702

$ws.Range("H136").Value = 6562.1304
$ws.Range("I136").Value = 7571.0625
$ws.Range("K136").Value = 22713.1875
$ws.Range("M136").Value = -20163.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 64690844
$ws.Range("I4").Value = 81495620
$ws.Range("J4").Value = 24359394
$ws.Range("K4").Value = 244486860
$ws.Range("L4").Value = 73078182
$ws.Range("M4").Value = -244486748
$ws.Range("N4").Value = -73078406

$ws.Range("H5").Value = 56994.723
$ws.Range("I5").Value = 794.1429000000001
$ws.Range("J5").Value = 92758.73
$ws.Range("K5").Value = 2382.4287
$ws.Range("L5").Value = 278276.19
$ws.Range("M5").Value = -2270.4287
$ws.Range("N5").Value = -278500.19

$ws.Range("H37").Value = 51833.332
$ws.Range("J37").Value = 51833.332
$ws.Range("L37").Value = 155499.996
$ws.Range("N37").Value = -155723.996

$ws.Range("H49").Value = 90
$ws.Range("I49").Value = 90
$ws.Range("K49").Value = 270
$ws.Range("M49").Value = -114

$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("M69").Value = -5189

$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("M72").Value = -13944

$ws.Range("H92").Value = 169.75
$ws.Range("J92").Value = 169.66667
$ws.Range("L92").Value = 509.00001
$ws.Range("N92").Value = -3005.00001

$ws.Range("H109").Value = 1642.8
$ws.Range("I109").Value = 1642.8
$ws.Range("K109").Value = 4928.4
$ws.Range("M109").Value = -3888.4

$ws.Range("H117").Value = 3676
$ws.Range("I117").Value = 6298.3335
$ws.Range("J117").Value = 2364.8333
$ws.Range("K117").Value = 18895.0005
$ws.Range("L117").Value = 7094.499899999999
$ws.Range("M117").Value = -15453.0005
$ws.Range("N117").Value = -13978.4999

$ws.Range("H133").Value = 2030
$ws.Range("I133").Value = 2030
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6090
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -1030
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 56994.723
$ws.Range("I135").Value = 794.1429000000001
$ws.Range("J135").Value = 92758.73
$ws.Range("K135").Value = 7147.2861
$ws.Range("L135").Value = 834828.5699999999
$ws.Range("M135").Value = -4612.2861
$ws.Range("N135").Value = -839898.5699999999

$ws.Range("H137").Value = 3794.4
$ws.Range("I137").Value = 3229.4546
$ws.Range("J137").Value = 5348
$ws.Range("K137").Value = 9688.363799999999
$ws.Range("L137").Value = 16044
$ws.Range("M137").Value = -4588.363799999999
$ws.Range("N137").Value = -26244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 32923378
$ws.Range("I80").Value = 43714836
$ws.Range("K80").Value = 43714836
$ws.Range("M80").Value = -43713838

$ws.Range("H83").Value = 32923378
$ws.Range("I83").Value = 43714836
$ws.Range("K83").Value = 218574180
$ws.Range("M83").Value = -218569188

$ws.Range("H107").Value = 632
$ws.Range("I107").Value = 967
$ws.Range("J107").Value = 520.3333
$ws.Range("K107").Value = 967
$ws.Range("L107").Value = 520.3333
$ws.Range("M107").Value = 953
$ws.Range("N107").Value = -4360.3333

$ws.Range("H122").Value = 251170.81
$ws.Range("I122").Value = 427189.38
$ws.Range("J122").Value = 4744.8
$ws.Range("K122").Value = 1281568.14
$ws.Range("L122").Value = 14234.4
$ws.Range("M122").Value = -1279118.14
$ws.Range("N122").Value = -19134.4

$ws.Range("H126").Value = 4773487.5
$ws.Range("I126").Value = 3249809.5
$ws.Range("J126").Value = 6414371.5
$ws.Range("K126").Value = 9749428.5
$ws.Range("L126").Value = 19243114.5
$ws.Range("M126").Value = -9746958.5
$ws.Range("N126").Value = -19248054.5

$ws.Range("H132").Value = 2814.077
$ws.Range("I132").Value = 2811.889
$ws.Range("J132").Value = 2819
$ws.Range("K132").Value = 8435.667000000001
$ws.Range("L132").Value = 8457
$ws.Range("M132").Value = -5905.667000000001
$ws.Range("N132").Value = -13517

$ws.Range("H136").Value = 11072.581
$ws.Range("J136").Value = 11072.581
$ws.Range("L136").Value = 33217.743
$ws.Range("N136").Value = -38317.743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6613.52
$ws.Range("I132").Value = 6439.915
$ws.Range("K132").Value = 19319.745
$ws.Range("M132").Value = -16789.745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 4500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 4500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 4500
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -4780

$ws.Range("H81").Value = 9809074
$ws.Range("I81").Value = 18521418
$ws.Range("K81").Value = 37042836
$ws.Range("M81").Value = -37041775

$ws.Range("H84").Value = 9809074
$ws.Range("I84").Value = 18521418
$ws.Range("K84").Value = 185214180
$ws.Range("M84").Value = -185208876

$ws.Range("H100").Value = 851.53845
$ws.Range("I100").Value = 910.9545000000001
$ws.Range("J100").Value = 524.75
$ws.Range("K100").Value = 1821.909
$ws.Range("L100").Value = 1049.5
$ws.Range("M100").Value = -1280.909
$ws.Range("N100").Value = -2131.5

$ws.Range("H113").Value = 2181.2942
$ws.Range("I113").Value = 738.8
$ws.Range("K113").Value = 2216.4
$ws.Range("M113").Value = -46.39999999999964
